$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "73 x 47" + [char]11 + "  4    7" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "3|    |"
$t.Cell(1,2).Range.Text = "27 x 81" + [char]11 + "  8    1" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "7|    |"
$t.Cell(1,3).Range.Text = "75 x 53" + [char]11 + "  5    3" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "5|    |"
$t.Cell(2,1).Range.Text = "79 x 72" + [char]11 + "  7    2" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "9|    |"
$t.Cell(2,2).Range.Text = "40 x 14" + [char]11 + "  1    4" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "0|    |"
$t.Cell(2,3).Range.Text = "92 x 86" + [char]11 + "  8    6" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "2|    |"
$t.Cell(3,1).Range.Text = "56 x 20" + [char]11 + "  2    0" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "6|    |"
$t.Cell(3,2).Range.Text = "71 x 32" + [char]11 + "  3    2" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "1|    |"
$t.Cell(3,3).Range.Text = "82 x 31" + [char]11 + "  3    1" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "2|    |"
$t.Cell(4,1).Range.Text = "70 x 83" + [char]11 + "  8    3" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "0|    |"
$t.Cell(4,2).Range.Text = "31 x 30" + [char]11 + "  3    0" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "1|    |"
$t.Cell(4,3).Range.Text = "55 x 75" + [char]11 + "  7    5" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "5|    |"
$t.Cell(5,1).Range.Text = "49 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "9|    |"
$t.Cell(5,2).Range.Text = "22 x 73" + [char]11 + "  7    3" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "2|    |"
$t.Cell(5,3).Range.Text = "17 x 85" + [char]11 + "  8    5" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "7|    |"

Write-Output "Updated 15 cells"
